$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Client Id (A2)
$ws.Range("A2").Value = "FCLXw697"

# Candidate ID (B2) - numeric
$ws.Range("B2").Value = 23081037

# User Name (C2)
$ws.Range("C2").Value = "hwgktnx59"

# Exam Password (D2)
$ws.Range("D2").Value = "h8!7#QBp"

# First Name (F2)
$ws.Range("F2").Value = "RpFFQEZI"

# Last Name (G2)
$ws.Range("G2").Value = "Feuq"
